# Update cryptocurrency price/volume data per upstream GitHub Actions scrape refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.298.59"
$ws.Range("E2").Value = "  -0.53%  "
$ws.Range("D3").Value = "1.669.05"
$ws.Range("E3").Value = "  -0.55%  "
$ws.Range("E4").Value = "  +0.46%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.75"
$ws.Range("E5").Value = "  +0.93%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5247"
$ws.Range("E6").Value = "  -1.21%  "
$ws.Range("E7").Value = "  +0.42%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2675"
$ws.Range("E8").Value = "  -0.92%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06336"
$ws.Range("E9").Value = "  -1.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.05"
$ws.Range("E10").Value = "  -3.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07765"
$ws.Range("E11").Value = "  -0.74%  "
$ws.Range("D12").Value = "1.675.11"
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.453"
$ws.Range("E13").Value = "  -1.35%  "
$ws.Range("D14").Value = "1.893.98"
$ws.Range("E14").Value = "  -0.70%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5501"
$ws.Range("E15").Value = "  -1.15%  "
$ws.Range("D16").Value = "0.0₅8284"
$ws.Range("E16").Value = "  -0.44%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.14"
$ws.Range("E17").Value = "  -0.73%  "
$ws.Range("D18").Value = "26.326.60"
$ws.Range("E18").Value = "  -0.58%  "
$ws.Range("E19").Value = "  +0.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.680"
$ws.Range("E20").Value = "  -1.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "195.25"
$ws.Range("E21").Value = "  +0.89%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.17"
$ws.Range("E22").Value = "  -1.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.089"
$ws.Range("E23").Value = "  -3.97%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.008"
$ws.Range("E24").Value = "  +0.63%  "
$ws.Range("E25").Value = "  -1.55%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1241"
$ws.Range("E26").Value = "  -3.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.222"
$ws.Range("E27").Value = "  -2.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.22"
$ws.Range("E28").Value = "  -0.11%  "
$ws.Range("E29").Value = "  -1.58%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.06177"
$ws.Range("E30").Value = "  -1.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.282"
$ws.Range("E31").Value = "  +0.67%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.604"
$ws.Range("E32").Value = "  -0.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.305"
$ws.Range("E33").Value = "  -4.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.635"
$ws.Range("E34").Value = "  -2.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9743"
$ws.Range("E35").Value = "  -3.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.426"
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.790"
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5781"
$ws.Range("E38").Value = "  -5.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01612"
$ws.Range("E39").Value = "  -1.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.038"
$ws.Range("E40").Value = "  -1.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8599"
$ws.Range("E41").Value = "  -0.56%  "
$ws.Range("E42").Value = "  +0.45%  "
$ws.Range("D43").Value = "1.025.75"
$ws.Range("E43").Value = "  -5.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.33"
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").Value = "1.809.49"
$ws.Range("E45").Value = "  -0.75%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₈109"
$ws.Range("E46").Value = "  +5.72%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "57.86"
$ws.Range("E47").Value = "  +1.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.010"
$ws.Range("E48").Value = "  +0.97%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.091"
$ws.Range("E49").Value = "  -0.67%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.490"
$ws.Range("E50").Value = "  +0.97%  "
$ws.Range("E51").Value = "  -0.41%  "
